$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before row 327 (old rows 327-402 shift down to 329-404)
$ws.Rows.Item(327).Insert()
$ws.Rows.Item(327).Insert()

# New row 327 values
$ws.Range("A327").Value2 = 11
$ws.Range("B327").Value2 = "Vega Monumental Concepción"
$ws.Range("C327").Value2 = "Bíobío"
$ws.Range("D327").Value2 = 44964
$ws.Range("E327").Value2 = 8
$ws.Range("F327").Value2 = 100112017
$ws.Range("G327").Value2 = "Apio"
$ws.Range("H327").Value2 = "Americana (o)"
$ws.Range("I327").Value2 = "Primera"
$ws.Range("J327").Value2 = 100
$ws.Range("K327").Value2 = 8000
$ws.Range("L327").Value2 = 8500
$ws.Range("M327").Value2 = 8250
$ws.Range("N327").Value2 = "$/docena de matas"
$ws.Range("O327").Value2 = "Provincia de Limarí"
$ws.Range("P327").Value2 = 1375
$ws.Range("Q327").Value2 = 6
$ws.Range("R327").Value2 = "Hortaliza"

# New row 328 values
$ws.Range("A328").Value2 = 11
$ws.Range("B328").Value2 = "Vega Monumental Concepción"
$ws.Range("C328").Value2 = "Bíobío"
$ws.Range("D328").Value2 = 44964
$ws.Range("E328").Value2 = 8
$ws.Range("F328").Value2 = 100112017
$ws.Range("G328").Value2 = "Apio"
$ws.Range("H328").Value2 = "Americana (o)"
$ws.Range("I328").Value2 = "Segunda"
$ws.Range("J328").Value2 = 50
$ws.Range("K328").Value2 = 7000
$ws.Range("L328").Value2 = 7000
$ws.Range("M328").Value2 = 7000
$ws.Range("N328").Value2 = "$/docena de matas"
$ws.Range("O328").Value2 = "Provincia de Limarí"
$ws.Range("P328").Value2 = 1167
$ws.Range("Q328").Value2 = 6
$ws.Range("R328").Value2 = "Hortaliza"
